# "updates to the notebooks."
#  - Bump the "last updated" date shown via the Date placeholders (slide
#    master, every slide layout, and the notes master) from 1/27/2020 to
#    2/10/2020.
#  - Bump the spelled-out date on the title slide from "January 27, 2020" to
#    "February 10, 2020".
#  - Tweak the Python-version text on the "About Python" slide.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "1/27/2020") {
            $tr.Text = "2/10/2020"
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master date placeholder: the shape's TextRange isn't directly
# writable on this object, so go through the Header/Footer date API instead.
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "2/10/2020"

# Title slide (Slide 1): spelled-out date in the subtitle.
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$subtitleRange = $subtitle.TextFrame.TextRange
$datePara = $subtitleRange.Paragraphs(1, 1)
if ($datePara.Text.TrimEnd("`r") -eq "January 27, 2020") {
    $datePara.Text = "February 10, 2020"
}

# Slide 4 ("About Python"): update version text.
$s4 = $p.Slides.Item(4)
$content = $s4.Shapes.Item(2)
$contentRange = $content.TextFrame.TextRange

$para1 = $contentRange.Paragraphs(1, 1)
if ($para1.Text.TrimEnd("`r") -eq "Developed by Guido van Rossum in the early 90s. Current versions: Python 3.7.2 and Python 2.7.15") {
    $para1.Text = "Developed by Guido van Rossum in the early 90s. Current versions: Python 3.7.x and Python 2.7.x"
}

$para2 = $contentRange.Paragraphs(2, 1)
if ($para2.Text.TrimEnd("`r") -eq "Python 3.7 is van Rossum's (and the larger community's) attempt at fixing some core issues with python 2.x") {
    $para2.Text = "Python 3.7.x is van Rossum’s (and the larger community’s) attempt at fixing some core issues with python 2.x"
}
